$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '39.040.19'
$ws.Range('E2').Value = '  -4.70%  '
$ws.Range('D3').Value = '2.214.39'
$ws.Range('E3').Value = '  -7.33%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '296.11'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -5.77%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '80.40'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -9.50%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.505'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -5.10%  '
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.457'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -7.92%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0773'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -8.28%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '28.02'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -10.52%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '45.78'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -13.14%  '
$ws.Range('E13').Value = '  -1.71%  '
$ws.Range('D14').Value = '2.551.66'
$ws.Range('E14').Value = '  -7.54%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.11'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -7.84%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '13.92'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -8.96%  '
$ws.Range('D17').Value = '2.230.55'
$ws.Range('E17').Value = '  -6.23%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.710'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -7.82%  '
$ws.Range('D19').Value = '38.933.39'
$ws.Range('E19').Value = '  -4.89%  '
$ws.Range('D20').Value = '0.0₃0858'
$ws.Range('E20').Value = '  -6.39%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.69'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -8.11%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '64.69'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -6.79%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.79'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -10.13%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '224.87'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -3.73%  '
$ws.Range('E25').Value = '  +0.00%  '
$ws.Range('E26').Value = '  -10.57%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.75'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -4.42%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '22.26'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -7.54%  '
$ws.Range('E29').Value = '  -2.33%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '8.89'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -5.37%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '148.49'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -3.54%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '31.13'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -8.93%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.999'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.21%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.74'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -9.36%  '
$ws.Range('E35').Value = '  -4.71%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0687'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -6.90%  '
$ws.Range('E37').Value = '  -4.71%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0958'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -4.59%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.62'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -7.15%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.59'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -8.36%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '14.32'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -11.22%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.61'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -6.54%  '
$ws.Range('D43').Value = '1.896.42'
$ws.Range('E43').Value = '  -4.04%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.07'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -12.17%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0254'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -6.74%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '16.06'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -9.52%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '8.80'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -9.01%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.51'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -9.74%  '
$ws.Range('D49').Value = '2.421.04'
$ws.Range('E49').Value = '  -7.85%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '69.80'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -4.34%  '
$ws.Range('E51').Value = '  -0.23%  '
